# Exercise 7 (R and Word)
# Insert the Exercise 7 write-up paragraphs right after the "Exercise 7"
# heading, before the trailing empty paragraph that precedes the sectPr.

$d = $word.ActiveDocument

# Locate the "Exercise 7" heading paragraph.
$exerciseIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Exercise 7") {
        $exerciseIndex = $i
    }
}

# The paragraph right after the heading is the (currently empty) trailing
# paragraph; insert the new content immediately before it so it ends up
# directly below "Exercise 7" while the trailing empty paragraph is kept.
$anchor = $d.Paragraphs.Item($exerciseIndex + 1)
$r = $anchor.Range

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
      <w:r>
        <w:t>Using the function pairs it’s</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> possible to see a potential correlation between expend and bad, lawyers, employ and pop.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Thus a first attempt to calculate the linear model will be made considering these factors.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">After assessing the first model (with the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>qqnorm</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> and the graph between fitted and </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>residuas</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>) it is possible to see that the</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>variances for the different fitted values is</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> concentrated </w:t>
      </w:r>
      <w:r>
        <w:t>in a region around</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> smaller value of fitted expenses.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>2nd Iteration:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve"> Considering that in the previous iteration lawyers and employ reject the null hypothesis the 2nd iteration </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>will  consider</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> only these factors (bad will also be considered since it has the biggest estimated </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>coeficient</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t>).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve"> Additionally, now we will calculate the regression considering interaction between the variables.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Assessing the new regression parameters:</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">The </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>qqnorm</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> presents a curved shape with some points far from the line, the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>qqnorm</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> shows a concentration around certain fitted values (&lt;1000).</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Given that bad is the variable with highest </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>coeficient</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> we will try to have better results by elevating bad to the power of 2.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">After a third iteration the </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>qqnorm</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> graph presents a better slope and distance between the points.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Also the residuals x </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>ffited</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> graph shows that the residuals are more spread.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">The model is: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>expenses</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> = 136.3 -8.186*bad -0.1297*lawyers + 0.08236*employ - 0.1440*bad^2</w:t>
      </w:r>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
Write-Output "Inserted Exercise 7 content after paragraph $exerciseIndex"
